# Adding emissions, adding storage start levels, a range of smaller
# improvements and fixes.
#
# The concrete change captured by this workbook edit is on the
# "Remove_units" sheet: the three rows describing the FR00 / Nuclear unit
# (one per scenario block: Distributed Energy 2040, National Trends 2025,
# Distributed Energy 2030) are removed, shifting the remaining rows up and
# shrinking the used range from A1:D23 to A1:D20.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("unitdata")
$ws2 = $wb.Worksheets.Item("Remove_units")

# Delete the FR00 / Nuclear rows (original row numbers 9, 16 and 23).
# Work bottom-to-top so the row indices for the rows still to be removed
# stay valid while earlier deletions shift everything below them up.
$ws2.Rows.Item(23).Delete()
$ws2.Rows.Item(16).Delete()
$ws2.Rows.Item(9).Delete()

# Restore the saved view/selection state recorded in the workbook:
#  - "unitdata" scrolls back to the top (no frozen topLeftCell) with its
#    remembered selection at E2,
#  - "Remove_units" stays the active tab, with its remembered selection
#    moved to B26.
$ws1.Range("E2").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("B26").Select() | Out-Null
